# Apply updated "F" column (collected/number) values to the "展览" and
# "全部类型" worksheets, matching the regenerated gh-pages data dump.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition) sheet updates
$wsExhibit.Range("F5").Value  = 1966
$wsExhibit.Range("F6").Value  = 5591
$wsExhibit.Range("F13").Value = 4437
$wsExhibit.Range("F18").Value = 42
$wsExhibit.Range("F19").Value = 42
$wsExhibit.Range("F25").Value = 11
$wsExhibit.Range("F33").Value = 337
$wsExhibit.Range("F42").Value = 335
$wsExhibit.Range("F45").Value = 19
$wsExhibit.Range("F46").Value = 420
$wsExhibit.Range("F48").Value = 222
$wsExhibit.Range("F49").Value = 143

# 全部类型 (All Types) sheet updates
$wsAll.Range("F6").Value  = 1966
$wsAll.Range("F7").Value  = 5591
$wsAll.Range("F13").Value = 4437
$wsAll.Range("F17").Value = 42
$wsAll.Range("F20").Value = 42
$wsAll.Range("F43").Value = 335
$wsAll.Range("F45").Value = 420
$wsAll.Range("F47").Value = 222
$wsAll.Range("F48").Value = 143
